# Add a new test case (Test case 5 / TC_05) to the QA sheet, following the
# same layout as the existing test case blocks (rows 25-32 "Test case 4"
# are used as the formatting template for the new rows 33-39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting from the "Test case 4" block (rows 25-32) down to the
# new "Test case 5" block (rows 33-39/40) ---------------------------------

$ws.Range("A25:E25").Copy() | Out-Null
$ws.Range("A33:E33").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A26:E26").Copy() | Out-Null
$ws.Range("A34:E34").PasteSpecial(-4122) | Out-Null

$ws.Range("A27:E27").Copy() | Out-Null
$ws.Range("A35:E35").PasteSpecial(-4122) | Out-Null

$ws.Range("A28:E28").Copy() | Out-Null
$ws.Range("A36:E36").PasteSpecial(-4122) | Out-Null

$ws.Range("A29:E29").Copy() | Out-Null
$ws.Range("A37:E37").PasteSpecial(-4122) | Out-Null

$ws.Range("A30:E30").Copy() | Out-Null
$ws.Range("A38:E38").PasteSpecial(-4122) | Out-Null

$ws.Range("A31:E31").Copy() | Out-Null
$ws.Range("A39:E39").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Title row (row 33) keeps the big 23.25pt font coming from the copied
# format, but its height needs to match the other title rows.
$ws.Rows(33).RowHeight = 23.25

# --- Fill in the text for the new test case ------------------------------

$ws.Range("A33").Value = "Test case 5"

$ws.Range("A35").Value = "ID"
$ws.Range("B35").Value = "TC_05"

$ws.Range("A36").Value = "Name"
$ws.Range("B36").Value = "Manual Testing"

$ws.Range("A38").Value = "№"
$ws.Range("B38").Value = "Description"
$ws.Range("C38").Value = "Expectations"
$ws.Range("D38").Value = "Result"
$ws.Range("E38").Value = "Status"

$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Testing the digital will form"
$ws.Range("C39").Value = "Things can be added"
$ws.Range("D39").Value = "Things can be added"
$ws.Range("E39").Value = "Pass"

# --- Merge the new title row the same way as the other title rows --------
$ws.Range("A33:E33").Merge() | Out-Null

# --- Update the active selection / view ----------------------------------
$ws.Range("H37").Select()
